$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-30 Sunday", "2024-07-01 Monday"),
    @("663÷9=", "206÷7="),
    @("283÷4=", "206÷8="),
    @("933÷2=", "821÷6="),
    @("703÷7=", "340÷7="),
    @("964÷8=", "394÷6="),
    @("918÷7=", "177÷5="),
    @("707÷2=", "570÷9="),
    @("261÷9=", "265÷4="),
    @("722÷9=", "663÷2="),
    @("935÷5=", "168÷7="),
    @("969÷8=", "775÷5="),
    @("125÷3=", "706÷6="),
    @("556÷6=", "279÷6="),
    @("958÷7=", "519÷2="),
    @("773÷4=", "984÷5="),
    @("423÷8=", "120÷2="),
    @("295÷3=", "635÷6="),
    @("839÷8=", "782÷2="),
    @("113÷3=", "807÷6="),
    @("606÷4=", "929÷5="),
    @("892÷2=", "120÷6="),
    @("319÷8=", "101÷8="),
    @("307÷2=", "948÷2="),
    @("181÷7=", "933÷9="),
    @("184÷6=", "104÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
